$d = $word.ActiveDocument

# 1. Remove the "Switching to the weekend data..." bullet entirely; the
#    following "These are not earth-shattering conclusions..." bullet
#    (and everything after it) shifts up to take its place.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Switching to the weekend data*") {
        $target = $p
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# 2. Replace the closing "rain" bullet's text (including the trailing
#    parenthetical run and the lastRenderedPageBreak marker that lived on
#    its first run) with the new weekend-snowfall conclusion paragraph.
$rng = $d.Content
$rng.Find.Execute("We failed to find any significant patterns in ridership lost due to rain. (See how much time we have left at this point to comment on this)", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "We also found that weekend traffic saw large decreases for even minimal amounts of snow. This differs from what we saw for weekdays and suggests that weekend routines are less compulsory and are much more sensitive to snowfall.", `
    2)
